$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextValue $ws "D2" "26.322.47"
Set-TextValue $ws "E2" "  -2.19%  "

# Row 3
Set-TextValue $ws "D3" "1.793.88"
Set-TextValue $ws "E3" "  -2.01%  "

# Row 4
Set-TextValue $ws "D4" "1.006"
Set-TextValue $ws "E4" "  -0.18%  "

# Row 5
Set-TextValue $ws "D5" "1.005"
Set-TextValue $ws "E5" "  -0.22%  "

# Row 6
Set-TextValue $ws "D6" "307.22"
Set-TextValue $ws "E6" "  -1.22%  "

# Row 7
Set-TextValue $ws "D7" "0.4511"
Set-TextValue $ws "E7" "  -1.58%  "

# Row 8
Set-TextValue $ws "D8" "0.3590"
Set-TextValue $ws "E8" "  -2.82%  "

# Row 9
Set-TextValue $ws "D9" "45.90"
Set-TextValue $ws "E9" "  +0.05%  "

# Row 10
Set-TextValue $ws "D10" "0.07078"
Set-TextValue $ws "E10" "  -1.48%  "

# Row 11
Set-TextValue $ws "D11" "0.8857"
Set-TextValue $ws "E11" "  +0.85%  "

# Row 12
Set-TextValue $ws "D12" "0.07808"
Set-TextValue $ws "E12" "  -0.67%  "

# Row 13
Set-TextValue $ws "D13" "19.41"
Set-TextValue $ws "E13" "  -1.20%  "

# Row 14
Set-TextValue $ws "D14" "1.844.53"
Set-TextValue $ws "E14" "  +1.05%  "

# Row 15
Set-TextValue $ws "D15" "5.290"
Set-TextValue $ws "E15" "  -0.85%  "

# Row 16
Set-TextValue $ws "D16" "6.337"
Set-TextValue $ws "E16" "  -0.87%  "

# Row 17
Set-TextValue $ws "D17" "84.81"
Set-TextValue $ws "E17" "  -2.55%  "

# Row 18
Set-TextValue $ws "D18" "1.007"
Set-TextValue $ws "E18" "  -0.18%  "

# Row 19
Set-TextValue $ws "D19" "0.000008513"
Set-TextValue $ws "E19" "  -2.22%  "

# Row 20
Set-TextValue $ws "D20" "1.006"
Set-TextValue $ws "E20" "  -0.15%  "

# Row 21
Set-TextValue $ws "D21" "14.25"
Set-TextValue $ws "E21" "  -1.53%  "

# Row 22
Set-TextValue $ws "D22" "26.341.29"
Set-TextValue $ws "E22" "  -2.21%  "

# Row 23
Set-TextValue $ws "D23" "4.991"
Set-TextValue $ws "E23" "  -0.23%  "

# Row 24
Set-TextValue $ws "D24" "2.054.24"
Set-TextValue $ws "E24" "  -0.65%  "

# Row 25
Set-TextValue $ws "D25" "10.51"

# Row 26
Set-TextValue $ws "D26" "1.976"
Set-TextValue $ws "E26" "  +0.09%  "

# Row 27
Set-TextValue $ws "D27" "152.12"
Set-TextValue $ws "E27" "  +0.83%  "

# Row 28
Set-TextValue $ws "E28" "  -1.96%  "

# Row 29
Set-TextValue $ws "D29" "2.025"
Set-TextValue $ws "E29" "  +2.77%  "

# Row 30
Set-TextValue $ws "D30" "112.05"
Set-TextValue $ws "E30" "  -1.62%  "

# Row 31
Set-TextValue $ws "D31" "4.862"
Set-TextValue $ws "E31" "  -1.17%  "

# Row 32
Set-TextValue $ws "D32" "0.08682"
Set-TextValue $ws "E32" "  -1.35%  "

# Row 33
Set-TextValue $ws "D33" "3.083"
Set-TextValue $ws "E33" "  +1.30%  "

# Row 34
Set-TextValue $ws "B34" "RenderToken"
Set-TextValue $ws "C34" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D34" "2.737"
Set-TextValue $ws "E34" "  +6.95%  "

# Row 35
Set-TextValue $ws "B35" "Filecoin"
Set-TextValue $ws "C35" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D35" "4.447"
Set-TextValue $ws "E35" "  -0.70%  "

# Row 36
Set-TextValue $ws "D36" "0.7236"
Set-TextValue $ws "E36" "  -3.96%  "

# Row 37
Set-TextValue $ws "D37" "1.104"
Set-TextValue $ws "E37" "  -2.49%  "

# Row 38
Set-TextValue $ws "D38" "1.005"
Set-TextValue $ws "E38" "  -0.02%  "

# Row 39
Set-TextValue $ws "D39" "1.069"
Set-TextValue $ws "E39" "  -1.55%  "

# Row 40
Set-TextValue $ws "D40" "0.01926"
Set-TextValue $ws "E40" "  -0.34%  "

# Row 41
Set-TextValue $ws "D41" "0.05091"
Set-TextValue $ws "E41" "  -0.77%  "

# Row 42
Set-TextValue $ws "D42" "2.866"
Set-TextValue $ws "E42" "  -1.25%  "

# Row 43
Set-TextValue $ws "D43" "0.5084"
Set-TextValue $ws "E43" "  +2.13%  "

# Row 44
Set-TextValue $ws "D44" "6.865"
Set-TextValue $ws "E44" "  -1.33%  "

# Row 45
Set-TextValue $ws "E45" "  -5.17%  "

# Row 46
Set-TextValue $ws "D46" "7.992"
Set-TextValue $ws "E46" "  -3.78%  "

# Row 47
Set-TextValue $ws "D47" "1.006"
Set-TextValue $ws "E47" "  -0.18%  "

# Row 48
Set-TextValue $ws "E48" "  -1.08%  "

# Row 49
Set-TextValue $ws "D49" "100.92"
Set-TextValue $ws "E49" "  -1.12%  "

# Row 50
Set-TextValue $ws "D50" "9.834"
Set-TextValue $ws "E50" "  -3.17%  "

# Row 51
Set-TextValue $ws "D51" "1.577"
Set-TextValue $ws "E51" "  -2.19%  "
